$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "62.496.34"
$ws.Range("E2").Value = "  -2.12%  "
Set-TextValue "D3" "2.655.35"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "551.12"
$ws.Range("E5").Value = "  -2.78%  "
Set-TextValue "D6" "155.86"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  +0.24%  "
Set-TextValue "D8" "0.585"
$ws.Range("E8").Value = "  -2.11%  "
Set-TextValue "D9" "0.104"
$ws.Range("E9").Value = "  -4.97%  "
$ws.Range("E10").Value = "  -3.86%  "
Set-TextValue "D11" "5.47"
$ws.Range("E11").Value = "  -3.46%  "
Set-TextValue "D12" "0.362"
$ws.Range("E12").Value = "  -4.87%  "
Set-TextValue "D13" "3.141.66"
$ws.Range("E13").Value = "  -2.38%  "
Set-TextValue "D14" "25.80"
$ws.Range("E14").Value = "  -4.41%  "
Set-TextValue "D15" "62.456.12"
$ws.Range("E15").Value = "  -1.90%  "
Set-TextValue "D16" "0.0000144"
$ws.Range("E16").Value = "  -3.94%  "
Set-TextValue "D17" "2.671.64"
$ws.Range("E17").Value = "  -2.41%  "
Set-TextValue "D18" "11.73"
$ws.Range("E18").Value = "  -7.20%  "
$ws.Range("E19").Value = "  -4.32%  "
Set-TextValue "D20" "341.44"
$ws.Range("E20").Value = "  -3.98%  "
Set-TextValue "D21" "6.12"
$ws.Range("E21").Value = "  -7.47%  "
$ws.Range("E22").Value = "  +0.03%  "
Set-TextValue "D23" "0.503"
$ws.Range("E23").Value = "  -3.79%  "
Set-TextValue "D24" "62.80"
$ws.Range("E24").Value = "  -2.92%  "
Set-TextValue "D25" "0.168"
$ws.Range("E25").Value = "  -1.05%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  +0.03%  "
Set-TextValue "D27" "8.03"
$ws.Range("E27").Value = "  -4.59%  "
Set-TextValue "D28" "1.37"
$ws.Range("E28").Value = "  +2.22%  "
Set-TextValue "D29" "0.0₃0836"
$ws.Range("E29").Value = "  -8.62%  "
Set-TextValue "D30" "7.14"
$ws.Range("E30").Value = "  -0.71%  "
Set-TextValue "D31" "1.90"
$ws.Range("E31").Value = "  -4.80%  "
Set-TextValue "D32" "161.81"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("E33").Value = "  +0.09%  "
Set-TextValue "D34" "4.77"
$ws.Range("E34").Value = "  -3.74%  "
Set-TextValue "D35" "1.43"
$ws.Range("E35").Value = "  -3.62%  "
Set-TextValue "D36" "19.25"
$ws.Range("E36").Value = "  -4.35%  "
Set-TextValue "D37" "1.75"
$ws.Range("E37").Value = "  -3.77%  "
Set-TextValue "D38" "334.26"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("E39").Value = "  -4.06%  "
Set-TextValue "D40" "0.912"
$ws.Range("E40").Value = "  -6.99%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "38.21"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "3.94"
$ws.Range("E42").Value = "  -3.99%  "
Set-TextValue "D43" "20.62"
$ws.Range("E43").Value = "  -6.24%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D44" "0.999"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D45" "0.611"
$ws.Range("E45").Value = "  -2.64%  "
Set-TextValue "D46" "19.89"
$ws.Range("E46").Value = "  -6.16%  "
Set-TextValue "D47" "10.99"
$ws.Range("E47").Value = "  -0.83%  "
Set-TextValue "D48" "0.0549"
$ws.Range("E48").Value = "  -6.54%  "
Set-TextValue "D49" "128.30"
$ws.Range("E49").Value = "  -3.10%  "
Set-TextValue "D50" "0.0957"
$ws.Range("E50").Value = "  -4.55%  "
Set-TextValue "D51" "0.0237"
$ws.Range("E51").Value = "  -5.89%  "
